# Append a new "TODO" bullet item to the end of the log, and move the
# _GoBack bookmark (which Word drops at the most recent edit point) into
# it, landing between "TODO" and ": Update player score accordingly."
# — exactly where the author's cursor would have been after typing
# "TODO" and before continuing with the rest of the sentence.

$d = $word.ActiveDocument

# The log's final bullet ("ERROR: The tiles were not being placed
# correctly. ...") currently also holds the _GoBack bookmark at its end.
$lastPara = $d.Paragraphs.Last

# Insert a new paragraph right after it; Word carries over the
# ListParagraph style + numbering (numId 1) and the run formatting
# (sz/szCs 24) from the paragraph it split off from.
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

# Type the full new bullet text in one go.
$newPara.Range.InsertBefore("TODO: Update player score accordingly.")

# Relocate _GoBack: drop it from the old spot and re-add it between
# "TODO" and ": Update player score accordingly." in the new bullet.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$goBackPos = $newPara.Range.Start + "TODO".Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
